$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $origStyle = $c.Style
    $c.Value = "'" + $val
    $c.Style = $origStyle
}

Set-TextValue 2 4 "25.884.71"
Set-TextValue 2 5 "  +0.00%  "
Set-TextValue 3 4 "1.737.13"
Set-TextValue 3 5 "  -0.24%  "
Set-TextValue 4 4 "0.9991"
Set-TextValue 4 5 "  -0.11%  "
Set-TextValue 5 4 "243.57"
Set-TextValue 5 5 "  +5.59%  "
Set-TextValue 6 4 "0.9993"
Set-TextValue 6 5 "  -0.03%  "
Set-TextValue 7 4 "0.5194"
Set-TextValue 7 5 "  -1.18%  "
Set-TextValue 8 4 "0.2746"
Set-TextValue 8 5 "  -0.12%  "
Set-TextValue 9 4 "39.90"
Set-TextValue 9 5 "  +1.20%  "
Set-TextValue 10 4 "0.06167"
Set-TextValue 10 5 "  +0.41%  "
Set-TextValue 11 4 "1.741.43"
Set-TextValue 11 5 "  -0.01%  "
Set-TextValue 12 4 "0.07188"
Set-TextValue 12 5 "  +1.26%  "
Set-TextValue 13 4 "15.00"
Set-TextValue 13 5 "  -1.24%  "
Set-TextValue 14 4 "0.6440"
Set-TextValue 14 5 "  +0.32%  "
Set-TextValue 15 4 "4.619"
Set-TextValue 15 5 "  +1.96%  "
Set-TextValue 16 4 "77.36"
Set-TextValue 16 5 "  -0.17%  "
Set-TextValue 17 5 "  +0.04%  "
Set-TextValue 18 4 "0.9990"
Set-TextValue 18 5 "  -0.10%  "
Set-TextValue 19 4 "25.913.00"
Set-TextValue 19 5 "  +0.14%  "
Set-TextValue 20 5 "  +2.24%  "
Set-TextValue 21 4 "0.000006796"
Set-TextValue 21 5 "  +1.72%  "
Set-TextValue 22 4 "1.966.06"
Set-TextValue 22 5 "  +0.17%  "
Set-TextValue 23 5 "  -0.41%  "
Set-TextValue 24 4 "8.639"
Set-TextValue 24 5 "  -1.24%  "
Set-TextValue 25 4 "5.296"
Set-TextValue 25 5 "  +2.65%  "
Set-TextValue 26 4 "136.52"
Set-TextValue 26 5 "  -2.61%  "
Set-TextValue 27 4 "1.523"
Set-TextValue 27 5 "  +0.35%  "
Set-TextValue 28 4 "15.24"
Set-TextValue 28 5 "  +0.60%  "
Set-TextValue 29 4 "1.771"
Set-TextValue 29 5 "  -1.06%  "
Set-TextValue 30 4 "105.21"
Set-TextValue 30 5 "  +2.41%  "
Set-TextValue 31 4 "3.962"
Set-TextValue 31 5 "  +6.37%  "
Set-TextValue 32 4 "0.08241"
Set-TextValue 32 5 "  -0.77%  "
Set-TextValue 33 4 "3.649"
Set-TextValue 33 5 "  +3.57%  "
Set-TextValue 34 4 "0.04681"
Set-TextValue 36 5 "  +1.41%  "
Set-TextValue 37 4 "0.6197"
Set-TextValue 37 5 "  -0.20%  "
Set-TextValue 38 4 "2.689"
Set-TextValue 38 5 "  +0.22%  "
Set-TextValue 39 5 "  +0.67%  "
Set-TextValue 40 4 "1.927"
Set-TextValue 40 5 "  -0.53%  "
Set-TextValue 41 4 "0.9989"
Set-TextValue 41 5 "  -0.01%  "
Set-TextValue 42 4 "100.27"
Set-TextValue 42 5 "  +0.31%  "
Set-TextValue 43 4 "0.3854"
Set-TextValue 43 5 "  -0.21%  "
Set-TextValue 44 4 "0.7473"
Set-TextValue 44 5 "  +2.12%  "
Set-TextValue 45 4 "5.002"
Set-TextValue 45 5 "  -0.09%  "
Set-TextValue 47 4 "6.255"
Set-TextValue 47 5 "  +0.02%  "
Set-TextValue 48 4 "55.00"
Set-TextValue 48 5 "  +2.39%  "
Set-TextValue 49 4 "0.05215"
Set-TextValue 49 5 "  -2.30%  "
Set-TextValue 50 4 "30.61"
Set-TextValue 50 5 "  +1.61%  "
Set-TextValue 51 4 "7.530"
Set-TextValue 51 5 "  -1.29%  "
